$d = $word.ActiveDocument

function New-WordXmlWrapper($innerBodyXml) {
    return "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:r='http://schemas.openxmlformats.org/officeDocument/2006/relationships'><w:body>" + $innerBodyXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so paragraph indices processed
# later in the script stay valid.
# ---------------------------------------------------------------------------

# --- III. "What can we improve?" (paragraph 15) ---------------------------
# Drop the stray <w:rPr><w:b/></w:rPr> left on the paragraph mark (pPr),
# keep everything else identical; then add the new body paragraph after it.
$p15 = $d.Paragraphs(15)
$xml15 = New-WordXmlWrapper("<w:p><w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:contextualSpacing/></w:pPr><w:r w:rsidRPr=`"00B82741`"><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">III. </w:t></w:r><w:bookmarkStart w:id=`"4`" w:name=`"WhatCanImprove`"/><w:r w:rsidR=`"00F81327`"><w:rPr><w:b/></w:rPr><w:t>What can we improve?</w:t></w:r><w:bookmarkEnd w:id=`"4`"/></w:p>")
$p15.Range.InsertXML($xml15)

$p15 = $d.Paragraphs(15)
$p15.Range.InsertParagraphAfter()
$newP16 = $d.Paragraphs(16)
$xmlNew16 = New-WordXmlWrapper("<w:p><w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:contextualSpacing/></w:pPr><w:r><w:t>We need to work on having more timely feedback. Commit to meeting when agreed. Team needs to add story points in addition to difficulty rating.</w:t></w:r></w:p>")
$newP16.Range.InsertXML($xmlNew16)

# --- II. "What went wrong?" (paragraph 13) ---------------------------------
# Add the new body paragraph (with the relocated _GoBack bookmark) after it.
$p13 = $d.Paragraphs(13)
$p13.Range.InsertParagraphAfter()
$newP14 = $d.Paragraphs(14)
$xmlNew14 = New-WordXmlWrapper("<w:p><w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:contextualSpacing/></w:pPr><w:r><w:t>Sick scrum master failed to communicate with team about rescheduled meeting. Team</w:t></w:r><w:r><w:t xml:space=`"preserve`"> didn't check documentation requirements.</w:t></w:r><w:r><w:t xml:space=`"preserve`"> We didn't </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>actually use</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> Trello.</w:t></w:r><w:bookmarkStart w:id=`"3`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"3`"/></w:p>")
$newP14.Range.InsertXML($xmlNew14)

# --- I. "What went well?" (paragraph 11) -----------------------------------
# Drop the stray <w:rPr><w:b/></w:rPr> left on the paragraph mark (pPr),
# keep everything else identical; then add the new body paragraph after it.
$p11 = $d.Paragraphs(11)
$xml11 = New-WordXmlWrapper("<w:p><w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:contextualSpacing/></w:pPr><w:r w:rsidRPr=`"00B82741`"><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">I. </w:t></w:r><w:bookmarkStart w:id=`"2`" w:name=`"WhatWentWell`"/><w:r w:rsidR=`"00F81327`"><w:rPr><w:b/></w:rPr><w:t>What went well?</w:t></w:r><w:bookmarkEnd w:id=`"2`"/></w:p>")
$p11.Range.InsertXML($xml11)

$p11 = $d.Paragraphs(11)
$p11.Range.InsertParagraphAfter()
$newP12 = $d.Paragraphs(12)
$xmlNew12 = New-WordXmlWrapper("<w:p><w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:contextualSpacing/></w:pPr><w:r><w:t xml:space=`"preserve`">Our team communication in meetings work </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:t>really well</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t>. Things were sequenced and flowed well moving forward. We had good organization.</w:t></w:r></w:p>")
$newP12.Range.InsertXML($xmlNew12)

# --- III. "What could improve?" (paragraph 9) -------------------------------
# Remove the old _GoBack bookmark pair sitting at the end of this paragraph.
# (Deleting the bookmark object directly avoids disturbing the hyperlink run
# that shares this paragraph.)
$d.Bookmarks("_GoBack").Delete()

Write-Output "edit complete"
